# Scheduled data-refresh: updates Universalis market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ -> columns H-N)
# for specific leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

# ALC row 76: Warding Off Temptation (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5292.8667
$ws.Range("I76").Value = 3515.889
$ws.Range("J76").Value = 7958.3335
$ws.Range("K76").Value = 3515.889
$ws.Range("L76").Value = 7958.3335
$ws.Range("M76").Value = -3200.889
$ws.Range("N76").Value = -8588.333500000001

# ALC row 79: The Garden of Arcane Delights (L) (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5292.8667
$ws.Range("I79").Value = 3515.889
$ws.Range("J79").Value = 7958.3335
$ws.Range("K79").Value = 3515.889
$ws.Range("L79").Value = 7958.3335
$ws.Range("M79").Value = -2423.889
$ws.Range("N79").Value = -10142.3335

# ALC row 103: Let Loose the Juice (Leve Item ID 19909)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 734.3
$ws.Range("I103").Value = 578.5714
$ws.Range("J103").Value = 1097.6666
$ws.Range("K103").Value = 1735.7142
$ws.Range("L103").Value = 3292.9998
$ws.Range("M103").Value = -1149.7142
$ws.Range("N103").Value = -4464.9998

# ALC row 132: Fast-forwarding Flora (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 92901.58
$ws.Range("J132").Value = 4672.5
$ws.Range("L132").Value = 14017.5
$ws.Range("N132").Value = -19077.5

# ALC row 138: All-night Crafting (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2640.2144
$ws.Range("I138").Value = 1817
$ws.Range("J138").Value = 3545.75
$ws.Range("K138").Value = 5451
$ws.Range("L138").Value = 10637.25
$ws.Range("M138").Value = -311
$ws.Range("N138").Value = -20917.25

# ARM row 32: Ingot We Trust (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16936.334
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 16936.334
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 16936.334
$ws.Range("N32").Value = -17510.334
$ws.Range("M32").ClearContents()

# ARM row 51: Everybody Cut Footloose (Leve Item ID 3858)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# ARM row 54: Family Secrets (Leve Item ID 2817)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 16666.334
$ws.Range("I54").Value = 5000
$ws.Range("K54").Value = 5000
$ws.Range("M54").Value = -4231

# ARM row 122: Haste for High Durium (Leve Item ID 36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2393.0264
$ws.Range("I122").Value = 1214.5652
$ws.Range("J122").Value = 4200
$ws.Range("K122").Value = 3643.6956
$ws.Range("L122").Value = 12600
$ws.Range("M122").Value = -1193.6956
$ws.Range("N122").Value = -17500

# ARM row 128: Heading toward Bankruptcy (Leve Item ID 34570)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# BSM row 99: Meddle in Metal (Leve Item ID 19943)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2093.3
$ws.Range("I99").Value = 1519.1428
$ws.Range("J99").Value = 3433
$ws.Range("K99").Value = 1519.1428
$ws.Range("L99").Value = 3433
$ws.Range("M99").Value = -21.14280000000008
$ws.Range("N99").Value = -6429

# BSM row 104: Hammer and Sails (Leve Item ID 19571)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H104").Value = 76507.2
$ws.Range("J104").Value = 76507.2
$ws.Range("L104").Value = 76507.2
$ws.Range("N104").Value = -83495.2

# BSM row 107: The Gold Experience (Leve Item ID 27706)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1795.55
$ws.Range("I107").Value = 1790.0526
$ws.Range("K107").Value = 1790.0526
$ws.Range("M107").Value = 129.9474

# CRP row 16: Raise the Roof (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 473.2857
$ws.Range("I16").Value = 587.5
$ws.Range("J16").Value = 321
$ws.Range("K16").Value = 587.5
$ws.Range("L16").Value = 321
$ws.Range("M16").Value = -300.5
$ws.Range("N16").Value = -895

# CRP row 113: Patient Patients (Leve Item ID 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 473.2857
$ws.Range("I113").Value = 587.5
$ws.Range("J113").Value = 321
$ws.Range("K113").Value = 587.5
$ws.Range("L113").Value = 321
$ws.Range("M113").Value = 1582.5
$ws.Range("N113").Value = -4661

# CUL row 7: It's Always Sunny in Vylbrand (Leve Item ID 4728)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1312734
$ws.Range("I7").Value = 274.58334
$ws.Range("K7").Value = 823.7500200000001
$ws.Range("M7").Value = -711.7500200000001

# CUL row 14: Keep Your Powder Dry (Leve Item ID 12886)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 143647.92
$ws.Range("I14").Value = 143647.92
$ws.Range("K14").Value = 430943.76
$ws.Range("M14").Value = -430770.76

# CUL row 70: Persona non Gratin (Leve Item ID 12867)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3635.2727
$ws.Range("I70").Value = 2001
$ws.Range("K70").Value = 6003
$ws.Range("M70").Value = -5688

# CUL row 73: Recipe for Disaster (L) (Leve Item ID 12867)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 3635.2727
$ws.Range("I73").Value = 2001
$ws.Range("K73").Value = 6003
$ws.Range("M73").Value = -4911

# CUL row 113: Can't Eat Just One (Leve Item ID 27843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1609.5883
$ws.Range("J113").Value = 1427.9231
$ws.Range("L113").Value = 4283.7693
$ws.Range("N113").Value = -8623.7693

# GSM row 52: It's My Business to Know Things (Leve Item ID 4147)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 49998.332
$ws.Range("J52").Value = 49998.332
$ws.Range("L52").Value = 49998.332
$ws.Range("N52").Value = -50516.332

# GSM row 101: Best-laid Planispheres (Leve Item ID 18513)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 28999.666
$ws.Range("J101").Value = 28999.666
$ws.Range("L101").Value = 28999.666
$ws.Range("N101").Value = -35489.666

# GSM row 102: Put the Metal to the Peddle (Leve Item ID 36169)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2623.8928
$ws.Range("I102").Value = 1422.3684
$ws.Range("K102").Value = 1422.3684
$ws.Range("M102").Value = 199.6315999999999

# GSM row 113: Copious Crystal Cannons (Leve Item ID 27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2747.739
$ws.Range("I113").Value = 1947
$ws.Range("K113").Value = 1947
$ws.Range("M113").Value = 223

# GSM row 126: Gold Rush Order (Leve Item ID 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4553.4346
$ws.Range("I126").Value = 2454.4
$ws.Range("J126").Value = 6168.077
$ws.Range("K126").Value = 7363.200000000001
$ws.Range("L126").Value = 18504.231
$ws.Range("M126").Value = -4893.200000000001
$ws.Range("N126").Value = -23444.231

# LTW row 22: Skin off Their Backs (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 79451.92
$ws.Range("J22").Value = 3092.7
$ws.Range("L22").Value = 3092.7
$ws.Range("N22").Value = -3682.7

# LTW row 27: Fire and Hide (Leve Item ID 5277)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 79451.92
$ws.Range("J27").Value = 3092.7
$ws.Range("L27").Value = 3092.7
$ws.Range("N27").Value = -3306.7

# LTW row 40: Best Served Toad (Leve Item ID 36248)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3417
$ws.Range("I40").Value = 3136
$ws.Range("J40").Value = 4166.3335
$ws.Range("K40").Value = 3136
$ws.Range("L40").Value = 4166.3335
$ws.Range("M40").Value = -3000
$ws.Range("N40").Value = -4438.3335

# LTW row 100: Tiger in the Sack (Leve Item ID 19995)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6026.222
$ws.Range("I100").Value = 2281
$ws.Range("J100").Value = 19134.5
$ws.Range("K100").Value = 2281
$ws.Range("L100").Value = 19134.5
$ws.Range("M100").Value = -1740
$ws.Range("N100").Value = -20216.5

# LTW row 132: Tenets of Tanning (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1155669.1
$ws.Range("I132").Value = 1574094.2
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4722282.6
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -4719752.6
$ws.Range("N132").Value = -20060

# LTW row 133: The Perfect Accessory (Leve Item ID 41903)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 140156.5
$ws.Range("J133").Value = 140156.5
$ws.Range("L133").Value = 140156.5
$ws.Range("N133").Value = -145216.5

# WVR row 32: Piling It On (Leve Item ID 3066)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 37500
$ws.Range("I32").Value = 25000
$ws.Range("J32").Value = 50000
$ws.Range("K32").Value = 25000
$ws.Range("L32").Value = 50000
$ws.Range("M32").Value = -24683
$ws.Range("N32").Value = -50634

# WVR row 107: Flax Wax (Leve Item ID 27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2157.476
$ws.Range("I107").Value = 832.625
$ws.Range("J107").Value = 6397
$ws.Range("K107").Value = 2497.875
$ws.Range("L107").Value = 19191
$ws.Range("M107").Value = -577.875
$ws.Range("N107").Value = -23031

# WVR row 110: Suits You (Leve Item ID 25825)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 149803.33
$ws.Range("J110").Value = 149803.33
$ws.Range("L110").Value = 149803.33
$ws.Range("N110").Value = -157983.33

# WVR row 122: Heavy Armoire (Leve Item ID 36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3257.4167
$ws.Range("I122").Value = 2391.3635
$ws.Range("K122").Value = 7174.0905
$ws.Range("M122").Value = -4724.0905

# WVR row 126: A Polished Purchase (Leve Item ID 36210)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3255
$ws.Range("I126").Value = 2872.2307
$ws.Range("K126").Value = 8616.6921
$ws.Range("M126").Value = -6146.6921

# WVR row 132: Comfy Cabins (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1194978
$ws.Range("I132").Value = 1921731.8
$ws.Range("J132").Value = 5744.5454
$ws.Range("K132").Value = 5765195.4
$ws.Range("L132").Value = 17233.6362
$ws.Range("M132").Value = -5762665.4
$ws.Range("N132").Value = -22293.6362

# WVR row 136: Weaving the Envelope (Leve Item ID 44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 16573570
$ws.Range("I136").Value = 20061426
$ws.Range("J136").Value = 6249.75
$ws.Range("K136").Value = 60184278
$ws.Range("L136").Value = 18749.25
$ws.Range("M136").Value = -60181728
$ws.Range("N136").Value = -23849.25

